# Rental Listings & Applications automation
# - adds the "Listing Sydney Homes" rental-listing columns to Sheet1
# - adds a new Sheet2 that holds the full listing record (incl. the
#   fields that didn't make it onto Sheet1: Two Beds and Sofa /
#   Available_Date / Do not come late / Occupants_Count)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# Sheet1 - extra header columns (D1:K1) and partial data row (D2:G2)
# ---------------------------------------------------------------

$ws1.Range("D1").Value = "Title"
$ws1.Range("E1").Value = "Moving Cost"
$ws1.Range("F1").Value = "Description"
$ws1.Range("G1").Value = "Target Rent"
$ws1.Range("H1").Value = "Furnishing "

$ws1.Range("I1").Value = "Available_Date"
$ws1.Range("I1").WrapText = $true
$ws1.Range("J1").Value = "Ideal_Tenant"
$ws1.Range("J1").WrapText = $true
$ws1.Range("K1").Value = "Occupants_Count"
$ws1.Range("K1").WrapText = $true

$ws1.Range("D2").Value = "Listing Sydney Homes"
$ws1.Range("D2").VerticalAlignment = -4160
$ws1.Range("E2").Value = 2000
$ws1.Range("E2").VerticalAlignment = -4160
$ws1.Range("F2").Value = "Listing my property for rental"
$ws1.Range("F2").WrapText = $true
$ws1.Range("G2").Value = 4000

$ws1.Rows.Item(1).RowHeight = 30
$ws1.Rows.Item(2).RowHeight = 75

$ws1.Columns.Item(4).ColumnWidth = 19.59
$ws1.Columns.Item(5).ColumnWidth = 11.02
$ws1.Columns.Item(6).ColumnWidth = 10.31
$ws1.Columns.Item(7).ColumnWidth = 10.31
$ws1.Columns.Item(8).ColumnWidth = 9.88

# Sheet1 is no longer the tab shown when the workbook opens, and the
# view is scrolled over to column C with the whole sheet selected.
$ws1.Cells.Select()

# ---------------------------------------------------------------
# Sheet2 (new) - full listing record
# ---------------------------------------------------------------

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "Title"
$ws2.Range("B1").Value = "Moving Cost"
$ws2.Range("C1").Value = "Description"
$ws2.Range("D1").Value = "Target Rent"
$ws2.Range("E1").Value = "Furnishing "

$ws2.Range("F1").Value = "Available_Date"
$ws2.Range("F1").WrapText = $true
$ws2.Range("G1").Value = "Ideal_Tenant"
$ws2.Range("G1").WrapText = $true
$ws2.Range("H1").Value = "Occupants_Count"
$ws2.Range("H1").WrapText = $true

$ws2.Range("A2").Value = "Listing Sydney Homes"
$ws2.Range("A2").VerticalAlignment = -4160
$ws2.Range("B2").Value = 2000
$ws2.Range("B2").VerticalAlignment = -4160
$ws2.Range("C2").Value = "Listing my property for rental"
$ws2.Range("C2").WrapText = $true
$ws2.Range("D2").Value = 4000
$ws2.Range("E2").Value = "Two Beds and Sofa"
$ws2.Range("E2").WrapText = $true
$ws2.Range("F2").Value = 43343
$ws2.Range("F2").NumberFormat = "m/d/yyyy"
$ws2.Range("G2").Value = "Do not come late"
$ws2.Range("G2").WrapText = $true
$ws2.Range("H2").Value = 2

$ws2.Rows.Item(1).RowHeight = 30
$ws2.Rows.Item(2).RowHeight = 45

$ws2.Columns.Item(1).ColumnWidth = 19.59
$ws2.Columns.Item(2).ColumnWidth = 11.02
$ws2.Columns.Item(3).ColumnWidth = 10.31
$ws2.Columns.Item(4).ColumnWidth = 10.31
$ws2.Columns.Item(5).ColumnWidth = 9.88
$ws2.Columns.Item(6).ColumnWidth = 9.88
$ws2.Columns.Item(7).ColumnWidth = 12.31
$ws2.Columns.Item(8).ColumnWidth = 8.02

$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

$ws2.Range("A1").Select()
$ws2.Activate()
